# fix subset method definition, remove unneeded plotSpec() override
$wb = $excel.ActiveWorkbook

# Add the new "components" worksheet at the end of the workbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "components"

# Header row
$ws.Range("B1").Value = "as-is"
$ws.Range("C1").Value = "almost as-is"
$ws.Range("D1").Value = "implement"
$ws.Range("E1").Value = "not supported"
$ws.Range("F1").Value = "ionize"
$ws.Range("G1").Value = "done"

# Data rows: Col A label, then X marks in the relevant columns, optional H note
$rows = @(
    @{ A = "`$";              mark = "B,G" },
    @{ A = "[";               mark = "C,G" },
    @{ A = "[[";              mark = "B,G" },
    @{ A = "as.data.table";   mark = "B,G" },
    @{ A = "componentInfo";   mark = "B,G" },
    @{ A = "componentTable";  mark = "B,G" },
    @{ A = "consensus";       mark = "E" },
    @{ A = "filter";          mark = "C,G" },
    @{ A = "findFGroup";      mark = "B,G" },
    @{ A = "groupNames";      mark = "B,G" },
    @{ A = "initialize";      mark = "C,G" },
    @{ A = "length";          mark = "B,G" },
    @{ A = "names";           mark = "B,G" },
    @{ A = "plotEIC";         mark = "B,D,G"; H = "Seems enough, assuming we're not planning to merge components" },
    @{ A = "plotEICHash";     mark = "B,G" },
    @{ A = "plotSpec";        mark = "B,D,G"; H = "Seems enough, assuming we're not planning to merge components" },
    @{ A = "plotSpecHash";    mark = "B,G" },
    @{ A = "show";            mark = "C,G" }
)

$r = 2
foreach ($row in $rows) {
    $ws.Range("A$r").Value = $row.A
    foreach ($col in $row.mark.Split(",")) {
        $ws.Range("$col$r").Value = "X"
    }
    if ($row.ContainsKey("H")) {
        $ws.Range("H$r").Value = $row.H
    }
    $r++
}

$ws.Columns.Item(1).EntireColumn.AutoFit() | Out-Null

$ws.Range("H17").Select()

# Select the compounds sheet and update its selection / deselect tab
$compounds = $wb.Worksheets.Item("compounds")
$compounds.Range("B1:G1").Select()

# Make the new components sheet the active/visible tab
$ws.Activate()

$wb.Windows.Item(1).ScrollWorkbookTabs(1)
